{"js": "// Helper: replace the text of a Word.Range with new text as a single run,\n// avoiding leftover formatting artifacts from the original run(s).\nfunction replaceRangeText(range, newText) {\n  const collapsed = range.insertText(\"\", Word.InsertLocation.replace);\n  return collapsed.insertText(newText, Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Remove the old \"_GoBack\" bookmark that used to sit around\n//    \"Studio win-client solution\" (it is about to be re-created at a\n//    different spot below). Doing this first avoids any ambiguity\n//    about which same-named bookmark gets removed later.\n// ---------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) \"... CloudSetupSdkSyncSampleSupport and Sample-Live_sync Assembly ...\"\n//    -> \"... CloudSetupSdkSyncSampleSupport Assembly ...\"\n//    The removed words are replaced by a (now-empty/collapsed) \"_GoBack\"\n//    bookmark left at that spot.\n// ---------------------------------------------------------------------\nlet results = body.search(\" and Sample-Live_sync\", { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  const collapsed = target.insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n  collapsed.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Merge the three runs of \"Close the Visual \" + \"Studio win-client\n//    solution\" + \".\" into a single run:\n//    \"Close the Visual Studio win-client solution.\"\n// ---------------------------------------------------------------------\nresults = body.search(\"Close the Visual Studio win-client solution.\", { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  replaceRangeText(results.items[0], \"Close the Visual Studio win-client solution.\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) Merge the two runs of \"Check the\" + \" Sample-Live-Sync project\n//    references and make sure that the Cloud.dll reference is\n//    \\u201cspecific version\\u201d, and that it is the correct version.\"\n//    into a single run.\n// ---------------------------------------------------------------------\nconst checkSentence =\n  \"Check the Sample-Live-Sync project references and make sure that the \" +\n  \"Cloud.dll reference is \\u201cspecific version\\u201d, and that it is the correct version.\";\n\nresults = body.search(checkSentence, { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  replaceRangeText(results.items[0], checkSentence);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Remove the old \"_GoBack\" bookmark that used to sit around\n#    \"Studio win-client solution\" (it is about to be re-created at a\n#    different spot below). Doing this first avoids any ambiguity\n#    about which same-named bookmark is affected later.\n# ---------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ---------------------------------------------------------------------\n# 2) \"... CloudSetupSdkSyncSampleSupport and Sample-Live_sync Assembly ...\"\n#    -> \"... CloudSetupSdkSyncSampleSupport Assembly ...\"\n#    The removed words are replaced by a (now-empty/collapsed) \"_GoBack\"\n#    bookmark left at that spot.\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \" and Sample-Live_sync\"\n$found = $find.Execute()\nif ($found) {\n    $rng.Text = \"\"\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n\n# ---------------------------------------------------------------------\n# 3) Merge the three runs of \"Close the Visual \" + \"Studio win-client\n#    solution\" + \".\" into a single run:\n#    \"Close the Visual Studio win-client solution.\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Close the Visual Studio win-client solution.\"\n$found = $find.Execute()\nif ($found) {\n    $rng.Text = \"\"\n    $rng2 = $d.Range($rng.Start, $rng.Start)\n    $rng2.Text = \"Close the Visual Studio win-client solution.\"\n}\n\n# ---------------------------------------------------------------------\n# 4) Merge the two runs of \"Check the\" + \" Sample-Live-Sync project\n#    references and make sure that the Cloud.dll reference is\n#    \"specific version\", and that it is the correct version.\" into a\n#    single run.\n# ---------------------------------------------------------------------\n$checkSentence = \"Check the Sample-Live-Sync project references and make sure that the Cloud.dll reference is \" + [char]8220 + \"specific version\" + [char]8221 + \", and that it is the correct version.\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = $checkSentence\n$found = $find.Execute()\nif ($found) {\n    $rng.Text = \"\"\n    $rng2 = $d.Range($rng.Start, $rng.Start)\n    $rng2.Text = $checkSentence\n}\n"}
